$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws 'D2' '68.459.54'
Set-TextCell $ws 'E2' '  +1.01%  '
Set-TextCell $ws 'D3' '3.274.09'
Set-TextCell $ws 'E3' '  +0.30%  '
Set-TextCell $ws 'D4' '1.00'
Set-TextCell $ws 'E4' '  +0.08%  '
Set-TextCell $ws 'D5' '582.38'
Set-TextCell $ws 'E5' '  +0.32%  '
Set-TextCell $ws 'D6' '185.26'
Set-TextCell $ws 'E6' '  +1.91%  '
Set-TextCell $ws 'E7' '  -0.04%  '
Set-TextCell $ws 'D8' '0.598'
Set-TextCell $ws 'E8' '  -0.54%  '
Set-TextCell $ws 'E9' '  -0.05%  '
Set-TextCell $ws 'D10' '6.65'
Set-TextCell $ws 'E10' '  -1.17%  '
Set-TextCell $ws 'D11' '0.422'
Set-TextCell $ws 'E11' '  +1.60%  '
Set-TextCell $ws 'D12' '3.841.46'
Set-TextCell $ws 'E12' '  +0.44%  '
Set-TextCell $ws 'E13' '  -0.20%  '
Set-TextCell $ws 'D14' '28.44'
Set-TextCell $ws 'E14' '  +0.14%  '
Set-TextCell $ws 'D15' '68.477.91'
Set-TextCell $ws 'E15' '  +1.09%  '
Set-TextCell $ws 'D16' '0.0000172'
Set-TextCell $ws 'E16' '  +1.57%  '
Set-TextCell $ws 'D17' '3.290.85'
Set-TextCell $ws 'E17' '  +1.38%  '
Set-TextCell $ws 'D18' '5.85'
Set-TextCell $ws 'E18' '  +0.26%  '
Set-TextCell $ws 'D19' '13.62'
Set-TextCell $ws 'E19' '  +0.87%  '
Set-TextCell $ws 'D20' '391.69'
Set-TextCell $ws 'E20' '  +4.22%  '
Set-TextCell $ws 'D21' '7.75'
Set-TextCell $ws 'E21' '  +1.36%  '
Set-TextCell $ws 'D22' '71.85'
Set-TextCell $ws 'E22' '  +0.90%  '
Set-TextCell $ws 'D23' '1.00'
Set-TextCell $ws 'E23' '  +0.11%  '
Set-TextCell $ws 'D24' '0.521'
Set-TextCell $ws 'E24' '  +1.80%  '
Set-TextCell $ws 'E25' '  +0.27%  '
Set-TextCell $ws 'E26' '  +4.27%  '
Set-TextCell $ws 'D27' '9.68'
Set-TextCell $ws 'E27' '  +0.93%  '
Set-TextCell $ws 'D28' '0.999'
Set-TextCell $ws 'E28' '  -0.01%  '
Set-TextCell $ws 'E29' '  +0.14%  '
Set-TextCell $ws 'D30' '5.73'
Set-TextCell $ws 'E30' '  +1.39%  '
Set-TextCell $ws 'D31' '23.13'
Set-TextCell $ws 'E31' '  +1.76%  '
Set-TextCell $ws 'D32' '7.20'
Set-TextCell $ws 'E32' '  +4.22%  '
Set-TextCell $ws 'D33' '1.30'
Set-TextCell $ws 'E33' '  +2.46%  '
Set-TextCell $ws 'D34' '0.999'
Set-TextCell $ws 'E34' '  +0.05%  '
Set-TextCell $ws 'D35' '164.26'
Set-TextCell $ws 'E35' '  +0.71%  '
Set-TextCell $ws 'D36' '1.51'
Set-TextCell $ws 'E36' '  +0.13%  '
Set-TextCell $ws 'D37' '1.94'
Set-TextCell $ws 'E37' '  +5.26%  '
Set-TextCell $ws 'D38' '0.829'
Set-TextCell $ws 'E38' '  -2.65%  '
Set-TextCell $ws 'D39' '26.93'
Set-TextCell $ws 'E39' '  +0.73%  '
Set-TextCell $ws 'D40' '4.62'
Set-TextCell $ws 'E40' '  -0.36%  '
Set-TextCell $ws 'D41' '6.63'
Set-TextCell $ws 'E41' '  -2.80%  '
Set-TextCell $ws 'D42' '2.55'
Set-TextCell $ws 'E42' '  -2.27%  '
Set-TextCell $ws 'D43' '41.36'
Set-TextCell $ws 'E43' '  +1.17%  '
Set-TextCell $ws 'D44' '0.0690'
Set-TextCell $ws 'E44' '  +1.90%  '
Set-TextCell $ws 'D45' '2.651.67'
Set-TextCell $ws 'E45' '  -1.58%  '
Set-TextCell $ws 'B46' 'InjectiveProtocol'
Set-TextCell $ws 'C46' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws 'D46' '25.29'
Set-TextCell $ws 'E46' '  -1.03%  '
Set-TextCell $ws 'B47' 'Bittensor'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws 'D47' '342.73'
Set-TextCell $ws 'E47' '  -2.54%  '
Set-TextCell $ws 'D48' '0.0282'
Set-TextCell $ws 'E48' '  +0.76%  '
Set-TextCell $ws 'B49' 'Arweave'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextCell $ws 'D49' '32.01'
Set-TextCell $ws 'E49' '  +2.43%  '
Set-TextCell $ws 'B50' 'Cosmos'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D50' '6.36'
Set-TextCell $ws 'E50' '  +3.49%  '
Set-TextCell $ws 'B51' 'ONDO'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextCell $ws 'D51' '0.995'
Set-TextCell $ws 'E51' '  -0.49%  '
